$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 256.30768
$ws.Range("H2").Value = 223.72223
$ws.Range("L2").Value = 139
$ws.Range("M2").Value = -143.30768
$ws.Range("N2").Value = -365
$ws.Range("K2").Value = 256.30768
$ws.Range("J2").Value = 139
$ws.Range("I4").Value = 57.166668
$ws.Range("H4").Value = 57.166668
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 56.833332
$ws.Range("K4").Value = 57.166668
$ws.Range("J4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("M5").Value = 57.090908
$ws.Range("K5").Value = 57.909092
$ws.Range("H5").Value = 65.583336
$ws.Range("I5").Value = 57.909092
$ws.Range("H17").Value = 741652.4
$ws.Range("L17").Value = 2527542.75
$ws.Range("J17").Value = 842514.25
$ws.Range("N17").Value = -2527878.75
$ws.Range("I19").Value = 831.1818
$ws.Range("L19").Value = 1918.8
$ws.Range("N19").Value = -2268.8
$ws.Range("M19").Value = -656.1818
$ws.Range("K19").Value = 831.1818
$ws.Range("J19").Value = 1918.8
$ws.Range("H19").Value = 1349.0952
$ws.Range("L40").Value = 870
$ws.Range("M40").Value = -861.7893999999999
$ws.Range("K40").Value = 1036.7894
$ws.Range("J40").Value = 870
$ws.Range("N40").Value = -1220
$ws.Range("I40").Value = 1036.7894
$ws.Range("H40").Value = 1020.9048
$ws.Range("K43").Value = 1599.5
$ws.Range("J43").Value = 11802.4
$ws.Range("I43").Value = 1599.5
$ws.Range("H43").Value = 10101.917
$ws.Range("L43").Value = 11802.4
$ws.Range("N43").Value = -11940.4
$ws.Range("M43").Value = -1530.5
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -4905.8
$ws.Range("K62").Value = 5529.8
$ws.Range("J62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H62").Value = 7206.125
$ws.Range("I62").Value = 5529.8
$ws.Range("H65").Value = 7206.125
$ws.Range("I65").Value = 5529.8
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -24529
$ws.Range("K65").Value = 27649
$ws.Range("J65").Value = 10000
$ws.Range("N65").Value = -56240
$ws.Range("H69").Value = 9538.429
$ws.Range("L69").Value = 28615.287
$ws.Range("N69").Value = -30363.287
$ws.Range("J69").Value = 9538.429
$ws.Range("I70").Value = 2500
$ws.Range("H70").Value = 3316.6667
$ws.Range("L70").Value = 11175
$ws.Range("N70").Value = -11715
$ws.Range("M70").Value = -7230
$ws.Range("K70").Value = 7500
$ws.Range("J70").Value = 3725
$ws.Range("L72").Value = 85845.861
$ws.Range("J72").Value = 9538.429
$ws.Range("N72").Value = -94581.861
$ws.Range("H72").Value = 9538.429
$ws.Range("I73").Value = 2500
$ws.Range("L73").Value = 11175
$ws.Range("N73").Value = -13047
$ws.Range("M73").Value = -6564
$ws.Range("K73").Value = 7500
$ws.Range("J73").Value = 3725
$ws.Range("H73").Value = 3316.6667
$ws.Range("H80").Value = 2472.3333
$ws.Range("I80").Value = 2360
$ws.Range("L80").Value = 7838.25
$ws.Range("M80").Value = -6082
$ws.Range("K80").Value = 7080
$ws.Range("J80").Value = 2612.75
$ws.Range("N80").Value = -9834.25
$ws.Range("H83").Value = 2472.3333
$ws.Range("I83").Value = 2360
$ws.Range("L83").Value = 23514.75
$ws.Range("N83").Value = -33498.75
$ws.Range("M83").Value = -16248
$ws.Range("K83").Value = 21240
$ws.Range("J83").Value = 2612.75
$ws.Range("H86").Value = 3747.5217
$ws.Range("I86").Value = 3094.5
$ws.Range("L86").Value = 5240.143
$ws.Range("M86").Value = -1971.5
$ws.Range("K86").Value = 3094.5
$ws.Range("J86").Value = 5240.143
$ws.Range("N86").Value = -7486.143
$ws.Range("I89").Value = 3094.5
$ws.Range("H89").Value = 3747.5217
$ws.Range("L89").Value = 26200.715
$ws.Range("N89").Value = -37432.715
$ws.Range("M89").Value = -9856.5
$ws.Range("K89").Value = 15472.5
$ws.Range("J89").Value = 5240.143
$ws.Range("L92").Value = 874.5
$ws.Range("N92").Value = -3370.5
$ws.Range("M92").Value = -7271.583000000001
$ws.Range("K92").Value = 8519.583000000001
$ws.Range("J92").Value = 874.5
$ws.Range("H92").Value = 7427.4287
$ws.Range("I92").Value = 8519.583000000001
$ws.Range("H99").Value = 5551.143
$ws.Range("H100").Value = 2484.8333
$ws.Range("L100").Value = 3450
$ws.Range("M100").Value = -1622.111
$ws.Range("K100").Value = 2163.111
$ws.Range("J100").Value = 3450
$ws.Range("N100").Value = -4532
$ws.Range("I100").Value = 2163.111
$ws.Range("M106").Value = -337
$ws.Range("K106").Value = 968
$ws.Range("I106").Value = 968
$ws.Range("H106").Value = 1501.091
$ws.Range("H112").Value = 1843.3334
$ws.Range("I112").Value = 1000
$ws.Range("L112").Value = 5846.25
$ws.Range("N112").Value = -8062.25
$ws.Range("M112").Value = -1892
$ws.Range("K112").Value = 3000
$ws.Range("J112").Value = 1948.75
$ws.Range("J113").Value = 4241
$ws.Range("I113").Value = 52779.094
$ws.Range("H113").Value = 35128.88
$ws.Range("L113").Value = 4241
$ws.Range("M113").Value = -49525.094
$ws.Range("N113").Value = -10749
$ws.Range("K113").Value = 52779.094
$ws.Range("I132").Value = 13931867
$ws.Range("H132").Value = 12898681
$ws.Range("L132").Value = 1501350
$ws.Range("N132").Value = -1506410
$ws.Range("M132").Value = -41793071
$ws.Range("K132").Value = 41795601
$ws.Range("J132").Value = 500450
$ws.Range("N137").Value = -13098.9999
$ws.Range("M137").Value = -1215.9231
$ws.Range("K137").Value = 3765.9231
$ws.Range("J137").Value = 2666.3333
$ws.Range("I137").Value = 1255.3077
$ws.Range("H137").Value = 1519.875
$ws.Range("L137").Value = 7998.999899999999
$ws.Range("K141").Value = 2985
$ws.Range("J141").Value = 2105
$ws.Range("I141").Value = 995
$ws.Range("H141").Value = 1550
$ws.Range("L141").Value = 6315
$ws.Range("N141").Value = -16675
$ws.Range("M141").Value = 2195

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M5").Value = -3580.6667
$ws.Range("K5").Value = 3692.6667
$ws.Range("H5").Value = 889.5
$ws.Range("I5").Value = 3692.6667
$ws.Range("I32").Value = 3079.5625
$ws.Range("H32").Value = 3484.1714
$ws.Range("M32").Value = -2792.5625
$ws.Range("K32").Value = 3079.5625
$ws.Range("N53").Value = -56364
$ws.Range("M53").Value = -9318
$ws.Range("K53").Value = 10000
$ws.Range("J53").Value = 55000
$ws.Range("I53").Value = 10000
$ws.Range("H53").Value = 25000
$ws.Range("L53").Value = 55000
$ws.Range("M61").Value = -1983.72
$ws.Range("K61").Value = 2195.72
$ws.Range("H61").Value = 2350.5312
$ws.Range("I61").Value = 2195.72
$ws.Range("M74").Value = -55706.223
$ws.Range("K74").Value = 56580.223
$ws.Range("H74").Value = 53449.57
$ws.Range("I74").Value = 56580.223
$ws.Range("L76").Value = 52999
$ws.Range("J76").Value = 52999
$ws.Range("N76").Value = -53675
$ws.Range("H76").Value = 52999
$ws.Range("I77").Value = 56580.223
$ws.Range("H77").Value = 53449.57
$ws.Range("M77").Value = -278533.115
$ws.Range("K77").Value = 282901.115
$ws.Range("J79").Value = 52999
$ws.Range("N79").Value = -55339
$ws.Range("H79").Value = 52999
$ws.Range("L79").Value = 52999
$ws.Range("M110").Value = -4732.4
$ws.Range("K110").Value = 6777.4
$ws.Range("H110").Value = 6598.788
$ws.Range("I110").Value = 6777.4
$ws.Range("I132").Value = 1984.119
$ws.Range("H132").Value = 2176.4082
$ws.Range("L132").Value = 9990.428400000001
$ws.Range("N132").Value = -15050.4284
$ws.Range("M132").Value = -3422.357
$ws.Range("K132").Value = 5952.357
$ws.Range("J132").Value = 3330.1428
$ws.Range("K136").Value = 6587.16
$ws.Range("I136").Value = 2195.72
$ws.Range("H136").Value = 2350.5312
$ws.Range("M136").Value = -4037.16

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 952.1539
$ws.Range("H54").Value = 12661.833
$ws.Range("I54").Value = 5194.2
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -50968
$ws.Range("M54").Value = -4710.2
$ws.Range("K54").Value = 5194.2
$ws.Range("J54").Value = 50000
$ws.Range("H86").Value = 17237.666
$ws.Range("I86").Value = 16872.445
$ws.Range("M86").Value = -15749.445
$ws.Range("K86").Value = 16872.445
$ws.Range("I89").Value = 16872.445
$ws.Range("H89").Value = 17237.666
$ws.Range("M89").Value = -78746.22500000001
$ws.Range("K89").Value = 84362.22500000001
$ws.Range("H99").Value = 2662.5833
$ws.Range("L99").Value = 11999
$ws.Range("M99").Value = -315.8181999999999
$ws.Range("K99").Value = 1813.8182
$ws.Range("J99").Value = 11999
$ws.Range("N99").Value = -14995
$ws.Range("I99").Value = 1813.8182
$ws.Range("N124").Value = -39819
$ws.Range("J124").Value = 29999
$ws.Range("H124").Value = 29999
$ws.Range("L124").Value = 29999
$ws.Range("H134").Value = 1930.3143
$ws.Range("I134").Value = 1554.1111
$ws.Range("L134").Value = 9600
$ws.Range("N134").Value = -14670
$ws.Range("M134").Value = -2127.3333
$ws.Range("K134").Value = 4662.3333
$ws.Range("J134").Value = 3200

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 5600
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5226
$ws.Range("J2").Value = 5000
$ws.Range("H31").Value = 40952.594
$ws.Range("I31").Value = 47337.863
$ws.Range("L31").Value = 12857.4
$ws.Range("M31").Value = -47042.863
$ws.Range("K31").Value = 47337.863
$ws.Range("J31").Value = 12857.4
$ws.Range("N31").Value = -13447.4
$ws.Range("H34").Value = 40952.594
$ws.Range("I34").Value = 47337.863
$ws.Range("L34").Value = 12857.4
$ws.Range("N34").Value = -13261.4
$ws.Range("M34").Value = -47135.863
$ws.Range("K34").Value = 47337.863
$ws.Range("J34").Value = 12857.4
$ws.Range("N57").ClearContents()
$ws.Range("J57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M58").Value = -1904.6924
$ws.Range("K58").Value = 2107.6924
$ws.Range("I58").Value = 2107.6924
$ws.Range("H58").Value = 2314.2856
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -5498.6665
$ws.Range("K62").Value = 6122.6665
$ws.Range("J62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H62").Value = 6010.4
$ws.Range("I62").Value = 6122.6665
$ws.Range("H65").Value = 6010.4
$ws.Range("I65").Value = 6122.6665
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -27493.3325
$ws.Range("K65").Value = 30613.3325
$ws.Range("J65").Value = 5000
$ws.Range("N65").Value = -31240
$ws.Range("H69").Value = 51118.625
$ws.Range("I69").Value = 45833.332
$ws.Range("M69").Value = -45084.332
$ws.Range("K69").Value = 45833.332
$ws.Range("M72").Value = -133755.996
$ws.Range("K72").Value = 137499.996
$ws.Range("I72").Value = 45833.332
$ws.Range("H72").Value = 51118.625
$ws.Range("H86").Value = 5265.6665
$ws.Range("I86").Value = 5265.6665
$ws.Range("M86").Value = -4142.6665
$ws.Range("K86").Value = 5265.6665
$ws.Range("I89").Value = 5265.6665
$ws.Range("H89").Value = 5265.6665
$ws.Range("M89").Value = -20712.3325
$ws.Range("K89").Value = 26328.3325
$ws.Range("H99").Value = 3731.3635
$ws.Range("L99").Value = 4143.75
$ws.Range("M99").Value = -2141.7222
$ws.Range("K99").Value = 3639.7222
$ws.Range("J99").Value = 4143.75
$ws.Range("N99").Value = -7139.75
$ws.Range("I99").Value = 3639.7222
$ws.Range("N105").Value = -4494
$ws.Range("M105").Value = 579.2106000000001
$ws.Range("K105").Value = 1167.7894
$ws.Range("J105").Value = 1000
$ws.Range("I105").Value = 1167.7894
$ws.Range("H105").Value = 1159.4
$ws.Range("L105").Value = 1000
$ws.Range("H107").Value = 939.94116
$ws.Range("I107").Value = 1102
$ws.Range("L107").Value = 642.8333
$ws.Range("M107").Value = 818
$ws.Range("K107").Value = 1102
$ws.Range("J107").Value = 642.8333
$ws.Range("N107").Value = -4482.8333
$ws.Range("J109").Value = 30321.5
$ws.Range("H109").Value = 30321.5
$ws.Range("L109").Value = 30321.5
$ws.Range("N109").Value = -32401.5
$ws.Range("I122").Value = 1028.0714
$ws.Range("H122").Value = 981.6875
$ws.Range("L122").Value = 1971
$ws.Range("M122").Value = -634.2142000000003
$ws.Range("K122").Value = 3084.2142
$ws.Range("J122").Value = 657
$ws.Range("N122").Value = -6871
$ws.Range("N126").Value = -17371.25
$ws.Range("M126").Value = -8449.1666
$ws.Range("K126").Value = 10919.1666
$ws.Range("J126").Value = 4143.75
$ws.Range("I126").Value = 3639.7222
$ws.Range("H126").Value = 3731.3635
$ws.Range("L126").Value = 12431.25
$ws.Range("I132").Value = 3197.074
$ws.Range("H132").Value = 3336.2903
$ws.Range("L132").Value = 12828
$ws.Range("N132").Value = -17888
$ws.Range("M132").Value = -7061.222
$ws.Range("K132").Value = 9591.222
$ws.Range("J132").Value = 4276
$ws.Range("J133").Value = 47587.5
$ws.Range("H133").Value = 47587.5
$ws.Range("L133").Value = 47587.5
$ws.Range("N133").Value = -52647.5
$ws.Range("H134").Value = 14677.944
$ws.Range("I134").Value = 8406.964
$ws.Range("M134").Value = -22685.892
$ws.Range("K134").Value = 25220.892
$ws.Range("K136").Value = 6323.0772
$ws.Range("I136").Value = 2107.6924
$ws.Range("H136").Value = 2314.2856
$ws.Range("M136").Value = -3773.0772
$ws.Range("H139").Value = 143900
$ws.Range("L139").Value = 143900
$ws.Range("N139").Value = -154180
$ws.Range("J139").Value = 143900

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 667900
$ws.Range("I14").Value = 667900
$ws.Range("M14").Value = -2003527
$ws.Range("K14").Value = 2003700
$ws.Range("I29").Value = 397.66666
$ws.Range("H29").Value = 677.8333
$ws.Range("L29").Value = 2874
$ws.Range("N29").Value = -3428
$ws.Range("M29").Value = -915.9999800000001
$ws.Range("K29").Value = 1192.99998
$ws.Range("J29").Value = 958
$ws.Range("H41").Value = 825.6
$ws.Range("I41").Value = 110.5
$ws.Range("L41").Value = 5207.1819
$ws.Range("N41").Value = -5883.1819
$ws.Range("M41").Value = 6.5
$ws.Range("K41").Value = 331.5
$ws.Range("J41").Value = 1735.7273
$ws.Range("N57").ClearContents()
$ws.Range("J57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M106").Value = -21384.5
$ws.Range("K106").Value = 22330.5
$ws.Range("I106").Value = 7443.5
$ws.Range("H106").Value = 6463.8
$ws.Range("H129").Value = 1496.1666
$ws.Range("L129").Value = 5843.571599999999
$ws.Range("J129").Value = 1947.8572
$ws.Range("N129").Value = -15843.5716
$ws.Range("K141").Value = 17691
$ws.Range("I141").Value = 5897
$ws.Range("H141").Value = 171586.33
$ws.Range("M141").Value = -12511

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 180.5
$ws.Range("H2").Value = 227.32
$ws.Range("M2").Value = -67.5
$ws.Range("K2").Value = 180.5
$ws.Range("H48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H49").Value = 32000
$ws.Range("L49").Value = 32000
$ws.Range("J49").Value = 32000
$ws.Range("N49").Value = -32368
$ws.Range("I70").Value = 9523.875
$ws.Range("H70").Value = 9303.904
$ws.Range("M70").Value = -9253.875
$ws.Range("K70").Value = 9523.875
$ws.Range("I73").Value = 9523.875
$ws.Range("M73").Value = -8587.875
$ws.Range("K73").Value = 9523.875
$ws.Range("H73").Value = 9303.904
$ws.Range("I97").Value = 1680.2
$ws.Range("H97").Value = 1507.6428
$ws.Range("L97").Value = 1076.25
$ws.Range("N97").Value = -2068.25
$ws.Range("M97").Value = -1184.2
$ws.Range("K97").Value = 1680.2
$ws.Range("J97").Value = 1076.25
$ws.Range("M102").Value = 90.36359999999991
$ws.Range("K102").Value = 1531.6364
$ws.Range("H102").Value = 1531.6364
$ws.Range("I102").Value = 1531.6364
$ws.Range("I122").Value = 1024
$ws.Range("H122").Value = 1299.2
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -622
$ws.Range("K122").Value = 3072
$ws.Range("J122").Value = 2400
$ws.Range("N122").Value = -12100
$ws.Range("N126").Value = -17007.8
$ws.Range("M126").Value = -73384.571
$ws.Range("K126").Value = 75854.571
$ws.Range("J126").Value = 4022.6
$ws.Range("I126").Value = 25284.857
$ws.Range("H126").Value = 19689.525
$ws.Range("L126").Value = 12067.8
$ws.Range("I132").Value = 235115.9
$ws.Range("H132").Value = 206632.34
$ws.Range("M132").Value = -702817.7
$ws.Range("K132").Value = 705347.7
$ws.Range("J136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("L136").Value = 0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 30668
$ws.Range("H7").Value = 25251
$ws.Range("M7").Value = -30556
$ws.Range("K7").Value = 30668
$ws.Range("H16").Value = 4607.1035
$ws.Range("I16").Value = 1646.5
$ws.Range("L16").Value = 7370.3335
$ws.Range("M16").Value = -1476.5
$ws.Range("K16").Value = 1646.5
$ws.Range("J16").Value = 7370.3335
$ws.Range("N16").Value = -7710.3335
$ws.Range("I22").Value = 1800
$ws.Range("H22").Value = 2093.625
$ws.Range("M22").Value = -1505
$ws.Range("K22").Value = 1800
$ws.Range("I27").Value = 1800
$ws.Range("H27").Value = 2093.625
$ws.Range("M27").Value = -1693
$ws.Range("K27").Value = 1800
$ws.Range("M40").Value = -4540
$ws.Range("K40").Value = 4676
$ws.Range("I40").Value = 4676
$ws.Range("H40").Value = 5340.8
$ws.Range("H41").Value = 32200
$ws.Range("L41").Value = 34400
$ws.Range("N41").Value = -35276
$ws.Range("J41").Value = 34400
$ws.Range("K43").Value = 24996
$ws.Range("I43").Value = 24996
$ws.Range("H43").Value = 24997.6
$ws.Range("M43").Value = -24803
$ws.Range("N46").Value = -1750
$ws.Range("J46").Value = 1374
$ws.Range("H46").Value = 1741.3334
$ws.Range("L46").Value = 1374
$ws.Range("J47").Value = 34400
$ws.Range("N47").Value = -35380
$ws.Range("H47").Value = 34400
$ws.Range("L47").Value = 34400
$ws.Range("J52").Value = 34400
$ws.Range("H52").Value = 34400
$ws.Range("L52").Value = 34400
$ws.Range("N52").Value = -34866
$ws.Range("I55").Value = 125.53333
$ws.Range("H55").Value = 142.75
$ws.Range("L55").Value = 194.4
$ws.Range("M55").Value = 47.46666999999999
$ws.Range("K55").Value = 125.53333
$ws.Range("J55").Value = 194.4
$ws.Range("N55").Value = -540.4
$ws.Range("M56").Value = -7309
$ws.Range("K56").Value = 8000
$ws.Range("I56").Value = 8000
$ws.Range("H56").Value = 8000
$ws.Range("M61").Value = -2311.5806
$ws.Range("K61").Value = 2513.5806
$ws.Range("H61").Value = 2791.4102
$ws.Range("I61").Value = 2513.5806
$ws.Range("H93").Value = 508266.72
$ws.Range("M93").Value = -926791.5600000001
$ws.Range("K93").Value = 928039.5600000001
$ws.Range("I93").Value = 928039.5600000001
$ws.Range("H100").Value = 3923.3333
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3382.3333
$ws.Range("K100").Value = 3923.3333
$ws.Range("J100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("I100").Value = 3923.3333
$ws.Range("N108").Value = -48305
$ws.Range("J108").Value = 40625
$ws.Range("H108").Value = 40625
$ws.Range("L108").Value = 40625
$ws.Range("I113").Value = 2513.5806
$ws.Range("H113").Value = 2791.4102
$ws.Range("M113").Value = -343.5805999999998
$ws.Range("K113").Value = 2513.5806
$ws.Range("I122").Value = 719540.6
$ws.Range("H122").Value = 594121.7
$ws.Range("L122").Value = 26499.999
$ws.Range("M122").Value = -2156171.8
$ws.Range("K122").Value = 2158621.8
$ws.Range("J122").Value = 8833.333000000001
$ws.Range("N122").Value = -31399.999
$ws.Range("M126").Value = -89534
$ws.Range("K126").Value = 92004
$ws.Range("I126").Value = 30668
$ws.Range("H126").Value = 25251
$ws.Range("H129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 33500
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
$ws.Range("J131").Value = 30000
$ws.Range("I132").Value = 3486.5217
$ws.Range("H132").Value = 3965.158
$ws.Range("L132").Value = 14097.201
$ws.Range("N132").Value = -19157.201
$ws.Range("M132").Value = -7929.5651
$ws.Range("K132").Value = 10459.5651
$ws.Range("J132").Value = 4699.067

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 22200
$ws.Range("L49").Value = 22200
$ws.Range("J49").Value = 22200
$ws.Range("N49").Value = -22660
$ws.Range("N57").ClearContents()
$ws.Range("J57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M64").Value = -11921
$ws.Range("K64").Value = 12169
$ws.Range("H64").Value = 12584.5
$ws.Range("I64").Value = 12169
$ws.Range("H67").Value = 12584.5
$ws.Range("M67").Value = -11311
$ws.Range("K67").Value = 12169
$ws.Range("I67").Value = 12169
$ws.Range("J109").Value = 37233.332
$ws.Range("H109").Value = 37233.332
$ws.Range("L109").Value = 37233.332
$ws.Range("N109").Value = -40007.332
$ws.Range("H122").Value = 4566.8
$ws.Range("L122").Value = 20998.5
$ws.Range("J122").Value = 6999.5
$ws.Range("N122").Value = -25898.5
$ws.Range("M126").Value = -18569
$ws.Range("K126").Value = 21039
$ws.Range("I126").Value = 7013
$ws.Range("H126").Value = 6797.6787
$ws.Range("I132").Value = 4475.375
$ws.Range("H132").Value = 4869.4614
$ws.Range("L132").Value = 16500
$ws.Range("N132").Value = -21560
$ws.Range("M132").Value = -10896.125
$ws.Range("K132").Value = 13426.125
$ws.Range("J132").Value = 5500
$ws.Range("K136").Value = 8906.25
$ws.Range("J136").Value = 2897.8333
$ws.Range("N136").Value = -13793.4999
$ws.Range("I136").Value = 2968.75
$ws.Range("H136").Value = 2938.3572
$ws.Range("L136").Value = 8693.499899999999
$ws.Range("M136").Value = -6356.25
